$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("396:397").Insert()
Write-Output "done"
